# Update Contig identifiers in the "Contig" column (column G) of Sheet1.
# Each mapping below gives the worksheet cells that hold a given old
# Contig value and the new value that should replace it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @(
    @{ Cells = @("G2", "G3", "G4"); Old = "Contig_2205_12.0141"; New = "Contig_2197_12.0141" },
    @{ Cells = @("G5");             Old = "Contig_2344_11.3427"; New = "Contig_2334_11.3427" },
    @{ Cells = @("G11");            Old = "Contig_2002_4.63142"; New = "Contig_2004_4.63142" },
    @{ Cells = @("G20");            Old = "Contig_176_36.506";   New = "Contig_172_36.506" },
    @{ Cells = @("G21");            Old = "Contig_189_39.2872";  New = "Contig_184_39.2872" },
    @{ Cells = @("G22");            Old = "Contig_45_77.3243";   New = "Contig_47_77.3243" },
    @{ Cells = @("G23");            Old = "Contig_65_57.9768";   New = "Contig_68_57.9768" },
    @{ Cells = @("G34");            Old = "Contig_36_117.344";   New = "Contig_38_117.344" },
    @{ Cells = @("G40");            Old = "Contig_178_57.0182";  New = "Contig_176_57.0182" }
)

foreach ($change in $changes) {
    foreach ($cellAddr in $change.Cells) {
        $cell = $ws.Range($cellAddr)
        if ($cell.Value2 -eq $change.Old) {
            $cell.Value = $change.New
        }
    }
}
